$wb = $excel.ActiveWorkbook

function Set-TextCell($cell, $text) {
    # Force the cell to be stored as literal text (matches the source
    # workbook's convention of inlineStr date/number-looking values)
    # instead of Excel's automatic number/date coercion.
    $cell.NumberFormat = "@"
    $cell.Value = $text
    $cell.ClearFormats()
}

# --- Sheet "HISP records": append rows 16-17 ---
$ws1 = $wb.Worksheets.Item("HISP records")

Set-TextCell $ws1.Cells.Item(16, 1) "02/03/2023"
$ws1.Cells.Item(16, 2).Value = "meeting internally on to discuss on EIR implementation plan"
Set-TextCell $ws1.Cells.Item(16, 3) "1"
$ws1.Cells.Item(16, 4).Value = "done"

Set-TextCell $ws1.Cells.Item(17, 1) "02/03/2023"
$ws1.Cells.Item(17, 2).Value = "wrote a reply letter on DHIS2 user accounts leakage notified by department of cyber security"
Set-TextCell $ws1.Cells.Item(17, 3) "1"
$ws1.Cells.Item(17, 4).Value = "done"

# --- Sheet "eBRS records": append rows 7-9 ---
$ws2 = $wb.Worksheets.Item("eBRS records")

Set-TextCell $ws2.Cells.Item(7, 1) "02/03/2023"
$ws2.Cells.Item(7, 2).Value = "meeting on pre-bid closing"
Set-TextCell $ws2.Cells.Item(7, 3) "2"
$ws2.Cells.Item(7, 4).Value = "done"

Set-TextCell $ws2.Cells.Item(8, 1) "05/03/2023"
$ws2.Cells.Item(8, 2).Value = "prepared presentation for bid closing to the committee"
Set-TextCell $ws2.Cells.Item(8, 3) "3"
$ws2.Cells.Item(8, 4).Value = "dne"

Set-TextCell $ws2.Cells.Item(9, 1) "05/03/2023"
$ws2.Cells.Item(9, 2).Value = "meeting with bid committee to finalize winning bidder of eBRS"
Set-TextCell $ws2.Cells.Item(9, 3) "2.5"
$ws2.Cells.Item(9, 4).Value = "done"
